$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: append a brand-new paragraph at the very end of the document
# body. Paragraphs.Add() keeps the new paragraph outside (after) the
# trailing bookmarkEnd and -- unlike Range.InsertParagraphAfter -- does
# not inherit stray character formatting from whatever text preceded
# the old end-of-body mark.
# ---------------------------------------------------------------------
function Add-EndParagraph([string]$style) {
    $endPos = $d.Content.End
    $insertAt = $d.Range($endPos, $endPos)
    $d.Paragraphs.Add($insertAt) | Out-Null
    $p = $d.Paragraphs.Last
    $p.Style = $style
    return $p
}

# =======================================================================
# Phase 1 -- create every new paragraph and drop in its plain text.
# (Character formatting is applied afterwards, in Phase 2: this runtime
# treats the most-recently-touched Font/Style property as "current
# typing formatting", so any format tweak made before a later
# Range.Text assignment can bleed into that later paragraph too.)
# =======================================================================

$pAbout = Add-EndParagraph "Heading2"
$pAbout.Range.Text = "About"

$pBody1 = Add-EndParagraph "FirstParagraph"
$body1Text = "This resume was written in Markdown, and styled with CSS, using the Markdown Resume project I built myself."
$pBody1.Range.Text = $body1Text
$b1 = $pBody1.Range.Start

$pBody2 = Add-EndParagraph "BodyText"
$body2Text = "With a simple script, I can generate a PDF, Word Document and HTML to embed into a website."
$pBody2.Range.Text = $body2Text
$b2 = $pBody2.Range.Start

$pLink = Add-EndParagraph "BlockText"
$linkStart = $pLink.Range.Start
$linkSentence = "Checkout my code for this project on GitHub."
$pLink.Range.Text = $linkSentence
$linkTextLen = $linkSentence.Length - 1   # exclude the trailing period

$pThanks = Add-EndParagraph "FirstParagraph"
$pThanks.Range.Text = "Thank you."

# =======================================================================
# Phase 2 -- apply character formatting / styles to the text just typed.
# =======================================================================

# "This resume was **written in Markdown**, and **styled with CSS**,
#  using the `Markdown Resume` project I built myself."
$d.Range($b1 + 0, $b1 + 11).Font.Italic = 1            # "This resume"
$d.Range($b1 + 16, $b1 + 35).Font.Bold = 1              # "written in Markdown"
$d.Range($b1 + 41, $b1 + 56).Font.Bold = 1              # "styled with CSS"
$d.Range($b1 + 68, $b1 + 83).Style = "VerbatimChar"     # "Markdown Resume"

# "With a simple `script`, I can generate a *PDF*, *Word Document* and
#  *HTML* to embed into a website."
$d.Range($b2 + 14, $b2 + 20).Style = "VerbatimChar"     # "script"
$d.Range($b2 + 39, $b2 + 42).Font.Italic = 1            # "PDF"
$d.Range($b2 + 44, $b2 + 57).Font.Italic = 1            # "Word Document"
$d.Range($b2 + 62, $b2 + 66).Font.Italic = 1            # "HTML"

# Hyperlink the GitHub sentence (but not the trailing period).
$linkRange = $d.Range($linkStart, $linkStart + $linkTextLen)
$d.Hyperlinks.Add($linkRange, "https://github.com/JasonToups/markdown-resume") | Out-Null

# =======================================================================
# Phase 3 -- wrap the whole new "About" block in a bookmark named
# "about", mirroring the existing "education" / "work-experience"
# section bookmarks. Positions are re-read fresh (not reused from
# Phase 1) because inserting the hyperlink above shifts the story's
# field-code bookkeeping, which would make any earlier-captured offsets
# stale.
# =======================================================================
$sectionStart = $pAbout.Range.Start
$sectionEnd = $pThanks.Range.End
$bookmarkRange = $d.Range($sectionStart, $sectionEnd)
$d.Bookmarks.Add("about", $bookmarkRange) | Out-Null

Write-Output "About section added"
